$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three remaining birth-year date values
$ws.Range("A2").Value = 36110
$ws.Range("A3").Value = 37179
$ws.Range("A4").Value = 28805

# Remove the now-unused rows 5 and 6 (shrinks the used range to A1:B4)
$ws.Rows("5:6").Delete()
